$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet1 "machinery reliability data": insert a new leading index column
#    ("python idx") with 0-based row numbers, matching the python-side index
#    used to initialise the ship model from this sheet.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Columns("A:A").Insert()
$ws1.Range("A1").Value = "python idx"
$ws1.Range("A1").WrapText = $true
$ws1.Range("A1").Font.Underline = $true
$ws1.Range("A1:A1").ColumnWidth = 12.833333333333332

for ($i = 0; $i -lt 13; $i++) {
    $ws1.Cells.Item($i + 2, 1).Value = $i
}

$ws1.Range("C37").Select()

# ---------------------------------------------------------------------------
# 2) Build the new "notes" sheet as a copy of "system structure data" (so it
#    keeps the same header/format), then adjust its rows: a new "test
#    system" entry replaces the old "fuel oil system" row, while the
#    "transmission system" / "diesel engine system" reference rows are kept.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2.Copy([System.Reflection.Missing]::Value, $lastSheet)
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "notes"

$ws3.Rows("2:2").Insert()
$ws3.Range("A2:B2").ClearFormats()
$ws3.Range("A2").Value = "test system"
$ws3.Range("B2").Value = "[(9,9), 4]"
$ws3.Rows("5:5").Delete()

$ws3.Range("A1:A1").ColumnWidth = 19.333333333333336
$ws3.Range("B1:B1").ColumnWidth = 29.666666666666668

$ws3.Range("A1:B4").Select()

# ---------------------------------------------------------------------------
# 3) Trim "system structure data" down to the single "fuel system" row
#    (renamed from "fuel oil system"); the other two reference rows now
#    live on the "notes" sheet instead.
# ---------------------------------------------------------------------------
$ws2.Rows("2:3").Delete()
$ws2.Range("A2").Value = "fuel system"

$ws2.Range("D10").Select()
